$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert the three new columns (shifts everything after them to the right)
#    - O:P  -> new "DA" / "TL Allowance" columns (inserted before old "Spcl
#              Allowance", which becomes Q)
#    - Y    -> new "Deducted allowance1" column (inserted before old
#              "total_deducations", which becomes Z)
# ---------------------------------------------------------------------------
$ws.Columns("O:P").Insert()
$ws.Columns("Y").Insert()

# ---------------------------------------------------------------------------
# 2) Header row (row 1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "SL #"
$ws.Range("B1").Value = "MONTH"
$ws.Range("C1").Value = "Emp. NAME"
$ws.Range("D1").Value = "DOJ"
$ws.Range("E1").Value = "STATUS"
$ws.Range("F1").Value = "DESIGNATION"
$ws.Range("G1").Value = "DEPARTMENT"
$ws.Range("H1").Value = "GROSS"
$ws.Range("I1").Value = "Per Month"
$ws.Range("J1").Value = "Actual Per Month"
$ws.Range("K1").Value = "Actual Days"
$ws.Range("L1").Value = "Working Days"
$ws.Range("M1").Value = "BASIC"
$ws.Range("N1").Value = "HRA"
$ws.Range("O1").Value = "DA"
$ws.Range("P1").Value = "TL Allowance"
$ws.Range("Q1").Value = "Spcl Allowance"
$ws.Range("R1").Value = "Arrears"
$ws.Range("S1").Value = "Gross Pay"
$ws.Range("T1").Value = "PF"
$ws.Range("U1").Value = "ESIC"
$ws.Range("V1").Value = "PT"
$ws.Range("W1").Value = "TDS"
$ws.Range("X1").Value = "Deductible Arrears"
$ws.Range("Y1").Value = "Deducted allowance1"
$ws.Range("Z1").Value = "total_deducations"
$ws.Range("AA1").Value = "NetPay"

# ---------------------------------------------------------------------------
# 3) Row 2 (existing employee, values updated)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "11-2014"
$ws.Range("C2").Value = "Priyanka Muddana"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "02/06/2014"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "Internship"
$ws.Range("F2").Value = "HR Manager"
$ws.Range("G2").Value = "HR"
$ws.Range("H2").Value = 750000.0
$ws.Range("I2").Value = 62500.0
$ws.Range("J2").Value = 58406.2
$ws.Range("K2").Value = 30.0
$ws.Range("L2").Value = 30.0
$ws.Range("M2").Value = 25000.0
$ws.Range("N2").Value = 6250.0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 27156.2
$ws.Range("R2").Value = 0.0
$ws.Range("S2").Value = 58406.2
$ws.Range("T2").Value = 3000.0
$ws.Range("U2").Value = 1022.11
$ws.Range("V2").Value = 0.0
$ws.Range("W2").Value = 13395.0
$ws.Range("X2").Value = 0.0
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 17417.1
$ws.Range("AA2").Value = 40989.1

# ---------------------------------------------------------------------------
# 4) Row 3 (new employee)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "11-2014"
$ws.Range("C3").Value = "Vidya Sagar pogiri"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "02/06/2014"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "Regular"
$ws.Range("F3").Value = "Junior Developer"
$ws.Range("G3").Value = "Development"
$ws.Range("H3").Value = 130000.0
$ws.Range("I3").Value = 10833.333333333334
$ws.Range("J3").Value = 10313.3
$ws.Range("K3").Value = 30.0
$ws.Range("L3").Value = 30.0
$ws.Range("M3").Value = 4333.33
$ws.Range("N3").Value = 1083.33
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 4896.67
$ws.Range("R3").Value = 0.0
$ws.Range("S3").Value = 10313.3
$ws.Range("T3").Value = 520.0
$ws.Range("U3").Value = 0.0
$ws.Range("V3").Value = 0.0
$ws.Range("W3").Value = 0.0
$ws.Range("X3").Value = 0.0
$ws.Range("Y3").Value = 0
$ws.Range("Z3").Value = 520.0
$ws.Range("AA3").Value = 9793.33

# ---------------------------------------------------------------------------
# 5) Column widths (A:AA) matching the re-computed "best fit" widths.
#    ColumnWidth is stored at a 1/7-character pixel grid internally, so the
#    inputs below are chosen as the closest achievable value to the target
#    OOXML <col width> figure.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 5.714285714285714
$ws.Columns("B").ColumnWidth = 9.0
$ws.Columns("C").ColumnWidth = 13.428571428571429
$ws.Columns("D").ColumnWidth = 13.428571428571429
$ws.Columns("E").ColumnWidth = 9.0
$ws.Columns("F").ColumnWidth = 14.428571428571429
$ws.Columns("G").ColumnWidth = 13.428571428571429
$ws.Columns("H").ColumnWidth = 10.142857142857142
$ws.Columns("I").ColumnWidth = 21.142857142857142
$ws.Columns("J").ColumnWidth = 13.428571428571429
$ws.Columns("K").ColumnWidth = 7.857142857142857
$ws.Columns("L").ColumnWidth = 10.142857142857142
$ws.Columns("M").ColumnWidth = 9.0
$ws.Columns("N").ColumnWidth = 9.0
$ws.Columns("O").ColumnWidth = 4.571428571428571
$ws.Columns("P").ColumnWidth = 10.142857142857142
$ws.Columns("Q").ColumnWidth = 11.142857142857142
$ws.Columns("R").ColumnWidth = 4.571428571428571
$ws.Columns("S").ColumnWidth = 9.0
$ws.Columns("T").ColumnWidth = 7.857142857142857
$ws.Columns("U").ColumnWidth = 9.0
$ws.Columns("V").ColumnWidth = 4.571428571428571
$ws.Columns("W").ColumnWidth = 9.0
$ws.Columns("X").ColumnWidth = 12.285714285714286
$ws.Columns("Y").ColumnWidth = 15.571428571428571
$ws.Columns("Z").ColumnWidth = 11.142857142857142
$ws.Columns("AA").ColumnWidth = 9.0
